# Update column G ("K") values per regenerated save_data (Strike# -> K)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 1
    3  = 1
    4  = 3
    5  = 0
    6  = 2
    7  = 1
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 1
    13 = 2
    14 = 0
    15 = 2
    17 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
